$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$A2 = @'
Sponsored
Visit Address By GS Thane - Consult an expert & visit site
raymonds-addressbygs.com
https://www.raymonds-addressbygs.com
Bookings Open Addres By GS Thane 6.1 Acre, 2/3/4 Bhk 1.30 Cr Ask Expert & Visit Site
'@
$ws.Range("A2").Value = $A2
$ws.Range("B2").Value = 'Home Bazaar Services Pvt Ltd'
$ws.Range("C2").Value = 'India'

$A3 = @'
Sponsored
The Address By GS Thane | 3, 4, 4.5 BHK Price ₹2.59Cr*
theaddressbygs-thane.in
https://www.theaddressbygs-thane.in
Launching The Address by GS at Pokhran Road Thane. Price Starts at ₹2.59 Cr*. Book Now.
Floor Plan & Pricing · Location · Project Highlights · Amenities Offered
'@
$ws.Range("A3").Value = $A3
$ws.Range("B3").Value = 'Prop Solutions4u'
$ws.Range("C3").Value = 'India'

$A4 = @'
Sponsored
New Raymond Tower in Thane | Presenting Homes in Thane west
raymondtenxera.com
https://www.raymondtenxera.com › official-site › brand
Prime Connectivity: School - 5 Mins | Metro - 3 Mins | Hospital - 2 Mins | Mall - 2 Mins. New...
Location Map · Configuration · Overview · Location Advantages · Contact Us · About Us
'@
$ws.Range("A4").Value = $A4
$ws.Range("B4").Value = 'Raymond Limited'
$ws.Range("C4").Value = 'India'

$A5 = @'
Sponsored
TenX Habitat Thane Launch | 2, 3 & 4 BHK Starts @ 1.41Cr*
homesfy-property.co.in
https://www.homesfy-property.co.in › tenx › thane
326,500 sq ft clubhouse | 2 and 3 BHK residences near Viviana Mall @ 1.41 CR*
'@
$ws.Range("A5").Value = $A5
$ws.Range("B5").Value = 'Homesfy Realty Limited'
$ws.Range("C5").Value = 'India'

$A6 = @'
Sponsored
RTMI Flats with Zero GST* | TenX Habitat by Raymond Realty
tenxhabitatraymondrealty.com
https://www.tenxhabitatraymondrealty.com
4,200 Sq.Ft. Multipurpose Hall | 28 Seater Mini Theatre | 2,400 Sq.Ft. Fully Equipped Gym. Book your 2 BHK Flat in Thane and get 50...
Real Estate Builders & Construction Company · Thane · Open ⋅ Closes 6 pm
Floor Plan · Floorplans Section · Location · Configurations Section · View Gallery · Amenities
'@
$ws.Range("A6").Value = $A6
$ws.Range("B6").Value = 'Raymond Limited'
$ws.Range("C6").Value = 'India'

# Remove rows 7-15 (shift cells up / delete entire rows)
$ws.Range("A7:C15").EntireRow.Delete() | Out-Null

